# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" between "2021-Q4" and "总计",
#    populated with the Q1-2022 fund holdings table (same column layout
#    as the "2021-Q4" sheet).
# 2. Update the "总计" (totals) summary sheet: insert a new first data
#    row for "2022-Q1" (6 funds, 0.07 亿元), pushing the existing
#    "2021-Q4" row down.

$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("2021-Q4")

# ---------------------------------------------------------------------
# 1. Create the "2022-Q1" worksheet, positioned right before "总计" so the
#    tab order becomes: 2021-Q4, 2022-Q1, 总计.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Add($wb.Worksheets.Item("总计"))
$ws.Name = "2022-Q1"

# NOTE: sheet handles in this host resolve by position, not stable
# identity - inserting "2022-Q1" *before* "总计" shifts 总计's index, so
# any reference captured beforehand (e.g. the argument just passed to
# Add()) now points at the newly inserted sheet instead. Re-resolve
# "总计" by name AFTER the insert so it refers to the right sheet.
$totalSheet = $wb.Worksheets.Item("总计")

# Header row (B1:H1) - copy formatting (bold/centered, style used by the
# "2021-Q4" sheet's header) from the source sheet, then set the text.
$src.Range("B1:H1").Copy()
$ws.Range("B1").PasteSpecial(-4122)

$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

# Column A (row index numbers 0..5) - copy the style used in the source
# sheet's A2 cell across the six data rows.
$src.Range("A2").Copy()
$ws.Range("A2:A7").PasteSpecial(-4122)

# Columns B:G mostly hold text-like values (fund code keeps leading
# zeros, the numeric-looking figures are stored as text) - force text
# formatting so Excel doesn't silently coerce them to numbers. G6:G7 are
# the exception (持有市值 rounds to plain number 0), so they are left at
# General format and set further below.
$ws.Range("B2:F7").NumberFormat = "@"
$ws.Range("G2:G5").NumberFormat = "@"

# Row 2 - 004258 国寿安保稳嘉混合A
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "004258"
$ws.Range("C2").Value = "国寿安保稳嘉混合A"
$ws.Range("D2").Value = "2.53"
$ws.Range("E2").Value = "22.03"
$ws.Range("F2").Value = "1.26"
$ws.Range("G2").Value = "0.0319"
$ws.Range("H2").Value = 5

# Row 3 - 005175 国寿安保消费新蓝海灵活配置混合
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "005175"
$ws.Range("C3").Value = "国寿安保消费新蓝海灵活配置混合"
$ws.Range("D3").Value = "0.71"
$ws.Range("E3").Value = "86.15"
$ws.Range("F3").Value = "3.29"
$ws.Range("G3").Value = "0.0234"
$ws.Range("H3").Value = 8

# Row 4 - 004301 国寿安保稳信混合A
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "004301"
$ws.Range("C4").Value = "国寿安保稳信混合A"
$ws.Range("D4").Value = "1.50"
$ws.Range("E4").Value = "20.03"
$ws.Range("F4").Value = "0.87"
$ws.Range("G4").Value = "0.0130"
$ws.Range("H4").Value = 9

# Row 5 - 004302 国寿安保稳信混合C
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "004302"
$ws.Range("C5").Value = "国寿安保稳信混合C"
$ws.Range("D5").Value = "0.01"
$ws.Range("E5").Value = "20.03"
$ws.Range("F5").Value = "0.87"
$ws.Range("G5").Value = "0.0001"
$ws.Range("H5").Value = 9

# Row 6 - 004259 国寿安保稳嘉混合C (持有市值 rounds to 0, stored as number 0)
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "004259"
$ws.Range("C6").Value = "国寿安保稳嘉混合C"
$ws.Range("D6").Value = "0.00"
$ws.Range("E6").Value = "22.03"
$ws.Range("F6").Value = "1.26"
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 5

# Row 7 - 015406 国寿安保稳信混合E (持有市值 rounds to 0, stored as number 0)
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "015406"
$ws.Range("C7").Value = "国寿安保稳信混合E"
$ws.Range("D7").Value = "0.00"
$ws.Range("E7").Value = "20.03"
$ws.Range("F7").Value = "0.87"
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 9

# The "@" text format was only needed transiently so Excel wouldn't
# coerce the numeric-looking strings above into real numbers; drop back
# to the workbook's default "Normal" style now that the values are
# locked in as text (keeps the cells unstyled, same as the source data).
$ws.Range("B2:F7").Style = "Normal"
$ws.Range("G2:G5").Style = "Normal"

# ---------------------------------------------------------------------
# 2. Update the "总计" sheet - add a "2022-Q1" row above the existing
#    "2021-Q4" row (new row 2, old row 2 becomes row 3).
# ---------------------------------------------------------------------

# Preserve the index-column styling (style used by A2) on the row the
# "2021-Q4" data is about to move into.
$totalSheet.Range("A2").Copy()
$totalSheet.Range("A3").PasteSpecial(-4122)

# Move the existing "2021-Q4" totals down to row 3.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2021-Q4"
$totalSheet.Range("C3").Value = 2
$totalSheet.Range("D3").Value = 0.17

# Write the new "2022-Q1" totals into row 2.
$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 6
$totalSheet.Range("D2").Value = 0.07
